$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The ID column ("3", "1", ...) looks numeric, so Excel would normally
# auto-type it as a real number. Temporarily mark those cells as text
# before assigning so they are stored the same way as the rest of the
# attendance sheet (plain text), then restore the default "Normal"
# style so no stray number-format styling is left behind.
$idCells = $ws.Range("A7:A8")
$idCells.NumberFormat = "@"

# Append two new attendance rows (7 and 8) after the existing data
$ws.Range("A7").Value = "3"
$ws.Range("B7").Value = "Jayaram"
$ws.Range("C7").Value = "24/11/2024"
$ws.Range("D7").Value = "13:29:33"

$ws.Range("A8").Value = "1"
$ws.Range("B8").Value = "Abhishek"
$ws.Range("C8").Value = "24/11/2024"
$ws.Range("D8").Value = "13:30:02"

$idCells.Style = "Normal"
